# Applies: "Added wire t connectors and weights to BOM"
#   - Insert "Wire T Connectors" row into the Power section (new row 7)
#   - Insert "Weights" row into the Hardware section (new row 20)
#   - Power Monitor (row 4) now carries an explicit unit price + formula total
#   - Final totals formula range grows to cover the two new rows
#   - Hyperlinks re-pointed at their (shifted) rows + two new hyperlinks added
#   - Selection marker moves to G5

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("BOM")

# ---------------------------------------------------------------------------
# 1. Insert the two new rows (bottom-most first so earlier indices are stable)
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).Insert()   # blank row that will become "Weights"
$ws.Rows.Item(7).Insert()    # blank row that will become "Wire T Connectors"

# ---------------------------------------------------------------------------
# 2. New row 7: Wire T Connectors (Power section)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Wire T Connectors"
$ws.Range("B7").Value = 14.99
$ws.Range("C7").Value = 1
$ws.Range("D7").Formula = "=C7*B7"
$ws.Range("B7:D7").NumberFormat = "$#,##0.00"
$ws.Range("C7").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 3. New row 20: Weights (Hardware section)
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Weights"
$ws.Range("B20").Value = 18.99
$ws.Range("C20").Value = 3
$ws.Range("D20").Formula = "=C20*B20"
$ws.Range("B20:D20").NumberFormat = "$#,##0.00"
$ws.Range("C20").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 4. Power Monitor (row 4) gains an explicit unit price + becomes a formula
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 13.88
$ws.Range("D4").Formula = "=C4*B4"

# ---------------------------------------------------------------------------
# 5. Fix up the grand-total SUM so it spans the newly inserted rows
# ---------------------------------------------------------------------------
$ws.Range("D25").Formula = "=SUM(D20:D24,D17:D18,D13:D15,D9:D11,D3:D7)"

# ---------------------------------------------------------------------------
# 6. Hyperlinks: engine does not auto-shift these on row insert, so clear and
#    re-create them all at their correct (post-insert) rows.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A3"),  "https://www.amazon.com/BOSYTRO-Switching-Transformer-Security-Industrial/dp/B0C591QLNR", "", "", "48V 25A 1200W DC Switching Power Supply") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"),  "https://www.amazon.com/dp/B0CZ36J9BY", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"),  "https://www.amazon.com/dp/B07FMTCHC1", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"),  "https://www.amazon.com/dp/B07FMRDP87", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"),  "https://www.aliexpress.us/item/3256805006838693.html", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"),  "https://shop.iflight.com/ipower-motor-gm5208-24-brushless-gimbal-motor-pro1347", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://www.aliexpress.us/item/1281133501.html", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A11"), "https://www.aliexpress.us/item/3256807861814061.html", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A13"), "https://odriverobotics.com/shop/odrive-s1", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A14"), "https://odriverobotics.com/shop/set-of-5-encoder-magnets", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A15"), "https://www.adafruit.com/product/4564?src=raspberrypi", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A17"), "https://www.amazon.com/dp/B0C2ZQHZ9T", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A18"), "https://www.amazon.com/Jameco-ValuePro-Stranded-Twisted-Hook-Up/dp/B0CT5SXZLM", "", "", "Twisted Pair") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A20"), "https://www.amazon.com/dp/B07V9FXWRV", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A21"), "https://www.aliexpress.us/item/3256804722090559.html", "", "", "Slider = MGN7C, 100mm") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A22"), "https://www.aliexpress.us/item/3256804722090559.html", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A23"), "https://www.amazon.com/dp/B0BZNST13Y", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A24"), "", "Hardware!A1", "", "Bolts, Nuts, and Pulleys") | Out-Null

# ---------------------------------------------------------------------------
# 7. Selection marker
# ---------------------------------------------------------------------------
$ws.Range("G5").Select()
